$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 93950
$ws.Range("J3").Value = 93950
$ws.Range("L3").Value = 93950
$ws.Range("N3").Value = -94178
$ws.Range("H43").Value = 18528.8
$ws.Range("I43").Value = 40266.332
$ws.Range("J43").Value = 9212.714
$ws.Range("K43").Value = 40266.332
$ws.Range("L43").Value = 9212.714
$ws.Range("M43").Value = -40197.332
$ws.Range("N43").Value = -9350.714
$ws.Range("H64").Value = 4987
$ws.Range("I64").Value = 2960.4736
$ws.Range("J64").Value = 9800
$ws.Range("K64").Value = 2960.4736
$ws.Range("L64").Value = 9800
$ws.Range("M64").Value = -2712.4736
$ws.Range("N64").Value = -10296
$ws.Range("H67").Value = 4987
$ws.Range("I67").Value = 2960.4736
$ws.Range("J67").Value = 9800
$ws.Range("K67").Value = 2960.4736
$ws.Range("L67").Value = 9800
$ws.Range("M67").Value = -2102.4736
$ws.Range("N67").Value = -11516
$ws.Range("H70").Value = 74548.78999999999
$ws.Range("I70").Value = 1733.3334
$ws.Range("J70").Value = 94407.55
$ws.Range("K70").Value = 5200.0002
$ws.Range("L70").Value = 283222.65
$ws.Range("M70").Value = -4930.0002
$ws.Range("N70").Value = -283762.65
$ws.Range("H73").Value = 74548.78999999999
$ws.Range("I73").Value = 1733.3334
$ws.Range("J73").Value = 94407.55
$ws.Range("K73").Value = 5200.0002
$ws.Range("L73").Value = 283222.65
$ws.Range("M73").Value = -4264.0002
$ws.Range("N73").Value = -285094.65
$ws.Range("H74").Value = 68073.12
$ws.Range("I74").Value = 103755.2
$ws.Range("K74").Value = 103755.2
$ws.Range("M74").Value = -102819.2
$ws.Range("H77").Value = 68073.12
$ws.Range("I77").Value = 103755.2
$ws.Range("K77").Value = 518776
$ws.Range("M77").Value = -514096
$ws.Range("H100").Value = 1365.35
$ws.Range("J100").Value = 413.8
$ws.Range("L100").Value = 413.8
$ws.Range("N100").Value = -1495.8
$ws.Range("H102").Value = 93950
$ws.Range("J102").Value = 93950
$ws.Range("L102").Value = 93950
$ws.Range("N102").Value = -100440
$ws.Range("H113").Value = 11513.294
$ws.Range("I113").Value = 7858.1113
$ws.Range("J113").Value = 15625.375
$ws.Range("K113").Value = 7858.1113
$ws.Range("L113").Value = 15625.375
$ws.Range("M113").Value = -4604.1113
$ws.Range("N113").Value = -22133.375
$ws.Range("H116").Value = 8580.799999999999
$ws.Range("I116").Value = 8726
$ws.Range("J116").Value = 8000
$ws.Range("K116").Value = 8726
$ws.Range("L116").Value = 8000
$ws.Range("M116").Value = -5284
$ws.Range("N116").Value = -14884
$ws.Range("H132").Value = 29503.666
$ws.Range("I132").Value = 1423.2258
$ws.Range("J132").Value = 203602.4
$ws.Range("K132").Value = 4269.6774
$ws.Range("L132").Value = 610807.2
$ws.Range("M132").Value = -1739.6774
$ws.Range("N132").Value = -615867.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10425.526
$ws.Range("I32").Value = 6766.3823
$ws.Range("K32").Value = 6766.3823
$ws.Range("M32").Value = -6479.3823
$ws.Range("H132").Value = 2898.658
$ws.Range("I132").Value = 2599.4194
$ws.Range("J132").Value = 4223.857
$ws.Range("K132").Value = 7798.2582
$ws.Range("L132").Value = 12671.571
$ws.Range("M132").Value = -5268.2582
$ws.Range("N132").Value = -17731.571

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2486.7144
$ws.Range("I86").Value = 2362.4
$ws.Range("J86").Value = 2797.5
$ws.Range("K86").Value = 2362.4
$ws.Range("L86").Value = 2797.5
$ws.Range("M86").Value = -1239.4
$ws.Range("N86").Value = -5043.5
$ws.Range("H89").Value = 2486.7144
$ws.Range("I89").Value = 2362.4
$ws.Range("J89").Value = 2797.5
$ws.Range("K89").Value = 11812
$ws.Range("L89").Value = 13987.5
$ws.Range("M89").Value = -6196
$ws.Range("N89").Value = -25219.5
$ws.Range("H134").Value = 1944.4736
$ws.Range("I134").Value = 1943.4259
$ws.Range("K134").Value = 5830.2777
$ws.Range("M134").Value = -3295.2777

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2842.8462
$ws.Range("I58").Value = 2543.625
$ws.Range("K58").Value = 2543.625
$ws.Range("M58").Value = -2340.625
$ws.Range("H62").Value = 6138.077
$ws.Range("J62").Value = 5415.143
$ws.Range("L62").Value = 5415.143
$ws.Range("N62").Value = -6663.143
$ws.Range("H65").Value = 6138.077
$ws.Range("J65").Value = 5415.143
$ws.Range("L65").Value = 27075.715
$ws.Range("N65").Value = -33315.715
$ws.Range("H132").Value = 3154.1943
$ws.Range("I132").Value = 2935.9656
$ws.Range("K132").Value = 8807.8968
$ws.Range("M132").Value = -6277.8968
$ws.Range("H136").Value = 2842.8462
$ws.Range("I136").Value = 2543.625
$ws.Range("K136").Value = 7630.875
$ws.Range("M136").Value = -5080.875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 8833
$ws.Range("J70").Value = 8799.6
$ws.Range("L70").Value = 26398.8
$ws.Range("N70").Value = -27028.8
$ws.Range("H73").Value = 8833
$ws.Range("J73").Value = 8799.6
$ws.Range("L73").Value = 26398.8
$ws.Range("N73").Value = -28582.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3000.75
$ws.Range("I80").Value = 2901.25
$ws.Range("J80").Value = 3199.75
$ws.Range("K80").Value = 2901.25
$ws.Range("L80").Value = 3199.75
$ws.Range("M80").Value = -1903.25
$ws.Range("N80").Value = -5195.75
$ws.Range("H83").Value = 3000.75
$ws.Range("I83").Value = 2901.25
$ws.Range("J83").Value = 3199.75
$ws.Range("K83").Value = 14506.25
$ws.Range("L83").Value = 15998.75
$ws.Range("M83").Value = -9514.25
$ws.Range("N83").Value = -25982.75
$ws.Range("H132").Value = 4226.1577
$ws.Range("I132").Value = 3275.6667
$ws.Range("J132").Value = 7790.5
$ws.Range("K132").Value = 9827.000100000001
$ws.Range("L132").Value = 23371.5
$ws.Range("M132").Value = -7297.000100000001
$ws.Range("N132").Value = -28431.5
$ws.Range("H141").Value = 88500
$ws.Range("I141").Value = 65000
$ws.Range("J141").Value = 96333.336
$ws.Range("K141").Value = 65000
$ws.Range("L141").Value = 96333.336
$ws.Range("M141").Value = -59820
$ws.Range("N141").Value = -106693.336

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2349.652
$ws.Range("I68").Value = 2347.3635
$ws.Range("J68").Value = 2400
$ws.Range("K68").Value = 2347.3635
$ws.Range("L68").Value = 2400
$ws.Range("M68").Value = -1598.3635
$ws.Range("N68").Value = -3898
$ws.Range("H71").Value = 2349.652
$ws.Range("I71").Value = 2347.3635
$ws.Range("J71").Value = 2400
$ws.Range("K71").Value = 11736.8175
$ws.Range("L71").Value = 12000
$ws.Range("M71").Value = -7992.817499999999
$ws.Range("N71").Value = -19488

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11107
$ws.Range("I62").Value = 7687.5
$ws.Range("K62").Value = 7687.5
$ws.Range("M62").Value = -7063.5
$ws.Range("H65").Value = 11107
$ws.Range("I65").Value = 7687.5
$ws.Range("K65").Value = 38437.5
$ws.Range("M65").Value = -35317.5
$ws.Range("H69").Value = 46499.5
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 46499.5
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 46499.5
$ws.Range("M69").ClearContents()  # was -44151
$ws.Range("N69").Value = -47997.5
$ws.Range("H72").Value = 46499.5
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 46499.5
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 139498.5
$ws.Range("M72").ClearContents()  # was -130956
$ws.Range("N72").Value = -146986.5

Write-Output "Applied 197 cell updates across 8 sheets"